# Box Plot Updates, Color Updates Main Figures
# Adjust jitter-plot marker / reference-line positions inside the single top-level
# group shape on slide 1. Target EMU offsets come from the authors diff; because the
# COM Left/Top properties round-trip through a single-precision (float32) point value
# before being re-expanded to EMU on save, the literal point values below were solved
# so that they serialize back to the *exact* target EMU (not just the nearest point).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)   # the single top-level group shape holding all figure elements

# pl6: (1397280,6003753) -> (1397280,6003451) EMU
$c = $grp.GroupItems.Item("pl6")
$c.Left = 110.02204895019531
$c.Top = 472.71270751953125
# pl7: (1397280,4306561) -> (1397280,4306263) EMU
$c = $grp.GroupItems.Item("pl7")
$c.Left = 110.02204895019531
$c.Top = 339.0758361816406
# pl8: (1397280,2609368) -> (1397280,2609075) EMU
$c = $grp.GroupItems.Item("pl8")
$c.Left = 110.02204895019531
$c.Top = 205.43898010253906
# pl9: (1397280,6852349) -> (1397280,6852045) EMU
$c = $grp.GroupItems.Item("pl9")
$c.Left = 110.02204895019531
$c.Top = 539.5311279296875
# pl10: (1397280,5155157) -> (1397280,5154857) EMU
$c = $grp.GroupItems.Item("pl10")
$c.Left = 110.02204895019531
$c.Top = 405.8942565917969
# pl11: (1397280,3457965) -> (1397280,3457669) EMU
$c = $grp.GroupItems.Item("pl11")
$c.Left = 110.02204895019531
$c.Top = 272.2574157714844
# pl12: (1397280,1760772) -> (1397280,1760481) EMU
$c = $grp.GroupItems.Item("pl12")
$c.Left = 110.02204895019531
$c.Top = 138.6205596923828
# pg15: (2423443,3299618) -> (2525642,3301105) EMU
$c = $grp.GroupItems.Item("pg15")
$c.Left = 198.8694610595703
$c.Top = 259.9295349121094
# pg16: (2253075,3299466) -> (1788073,3301740) EMU
$c = $grp.GroupItems.Item("pg16")
$c.Left = 140.79315185546875
$c.Top = 259.97955322265625
# pg17: (2318892,3893637) -> (2156969,3895131) EMU
$c = $grp.GroupItems.Item("pg17")
$c.Left = 169.840087890625
$c.Top = 306.7032470703125
# pg18: (2604626,3317827) -> (2494833,3317030) EMU
$c = $grp.GroupItems.Item("pg18")
$c.Left = 196.4435577392578
$c.Top = 261.1834716796875
# pg19: (2135399,3752547) -> (2015428,3753257) EMU
$c = $grp.GroupItems.Item("pg19")
$c.Left = 158.69512939453125
$c.Top = 295.5320739746094
# pg20: (1896327,3701760) -> (2540369,3701822) EMU
$c = $grp.GroupItems.Item("pg20")
$c.Left = 200.02906799316406
$c.Top = 291.4820556640625
# pl23: (4284320,6003753) -> (4284320,6003451) EMU
$c = $grp.GroupItems.Item("pl23")
$c.Left = 337.3480529785156
$c.Top = 472.71270751953125
# pl24: (4284320,4306561) -> (4284320,4306263) EMU
$c = $grp.GroupItems.Item("pl24")
$c.Left = 337.3480529785156
$c.Top = 339.0758361816406
# pl25: (4284320,2609368) -> (4284320,2609075) EMU
$c = $grp.GroupItems.Item("pl25")
$c.Left = 337.3480529785156
$c.Top = 205.43898010253906
# pl26: (4284320,6852349) -> (4284320,6852045) EMU
$c = $grp.GroupItems.Item("pl26")
$c.Left = 337.3480529785156
$c.Top = 539.5311279296875
# pl27: (4284320,5155157) -> (4284320,5154857) EMU
$c = $grp.GroupItems.Item("pl27")
$c.Left = 337.3480529785156
$c.Top = 405.8942565917969
# pl28: (4284320,3457965) -> (4284320,3457669) EMU
$c = $grp.GroupItems.Item("pl28")
$c.Left = 337.3480529785156
$c.Top = 272.2574157714844
# pl29: (4284320,1760772) -> (4284320,1760481) EMU
$c = $grp.GroupItems.Item("pl29")
$c.Left = 337.3480529785156
$c.Top = 138.6205596923828
# pg32: (5923627,6020770) -> (5834431,6020200) EMU
$c = $grp.GroupItems.Item("pg32")
$c.Left = 459.4040222167969
$c.Top = 474.0315246582031
# pg33: (6660827,5686827) -> (6149706,5686675) EMU
$c = $grp.GroupItems.Item("pg33")
$c.Left = 484.22882080078125
$c.Top = 447.76971435546875
# pl36: (7171360,6003753) -> (7171360,6003451) EMU
$c = $grp.GroupItems.Item("pl36")
$c.Left = 564.674072265625
$c.Top = 472.71270751953125
# pl37: (7171360,4306561) -> (7171360,4306263) EMU
$c = $grp.GroupItems.Item("pl37")
$c.Left = 564.674072265625
$c.Top = 339.0758361816406
# pl38: (7171360,2609368) -> (7171360,2609075) EMU
$c = $grp.GroupItems.Item("pl38")
$c.Left = 564.674072265625
$c.Top = 205.43898010253906
# pl39: (7171360,6852349) -> (7171360,6852045) EMU
$c = $grp.GroupItems.Item("pl39")
$c.Left = 564.674072265625
$c.Top = 539.5311279296875
# pl40: (7171360,5155157) -> (7171360,5154857) EMU
$c = $grp.GroupItems.Item("pl40")
$c.Left = 564.674072265625
$c.Top = 405.8942565917969
# pl41: (7171360,3457965) -> (7171360,3457669) EMU
$c = $grp.GroupItems.Item("pl41")
$c.Left = 564.674072265625
$c.Top = 272.2574157714844
# pl42: (7171360,1760772) -> (7171360,1760481) EMU
$c = $grp.GroupItems.Item("pl42")
$c.Left = 564.674072265625
$c.Top = 138.6205596923828
# pg45: (7453217,2978404) -> (7913146,2977284) EMU
$c = $grp.GroupItems.Item("pg45")
$c.Left = 623.0823974609375
$c.Top = 234.43182373046875
# pg46: (7853039,3246830) -> (7837029,3246101) EMU
$c = $grp.GroupItems.Item("pg46")
$c.Left = 617.0889282226562
$c.Top = 255.5985107421875
# pg47: (8680681,6634049) -> (9150355,6634049) EMU
$c = $grp.GroupItems.Item("pg47")
$c.Left = 720.5004272460938
$c.Top = 522.3660888671875
# pg48: (9380955,6463806) -> (9492295,6464172) EMU
$c = $grp.GroupItems.Item("pg48")
$c.Left = 747.4248046875
$c.Top = 508.98992919921875
# pg49: (9627766,5159595) -> (8867464,5160585) EMU
$c = $grp.GroupItems.Item("pg49")
$c.Left = 698.2255249023438
$c.Top = 406.34527587890625
# pg50: (7624975,3423975) -> (8230359,3421765) EMU
$c = $grp.GroupItems.Item("pg50")
$c.Left = 648.059814453125
$c.Top = 269.4303283691406
# pg51: (8833816,6144943) -> (9479043,6144594) EMU
$c = $grp.GroupItems.Item("pg51")
$c.Left = 746.38134765625
$c.Top = 483.8263244628906
# pg52: (8929301,5782988) -> (9504320,5781289) EMU
$c = $grp.GroupItems.Item("pg52")
$c.Left = 748.3717041015625
$c.Top = 455.2196350097656
# pg53: (7490425,4650846) -> (8144319,4651339) EMU
$c = $grp.GroupItems.Item("pg53")
$c.Left = 641.2849731445312
$c.Top = 366.2471923828125
# pg54: (8357919,2278102) -> (8208886,2279950) EMU
$c = $grp.GroupItems.Item("pg54")
$c.Left = 646.3690185546875
$c.Top = 179.5236358642578
# pg55: (8082146,2880590) -> (8205657,2878415) EMU
$c = $grp.GroupItems.Item("pg55")
$c.Left = 646.11474609375
$c.Top = 226.6468505859375
# pg56: (7851440,2465416) -> (8350502,2466290) EMU
$c = $grp.GroupItems.Item("pg56")
$c.Left = 657.5198974609375
$c.Top = 194.19607543945312
# pg57: (7753683,4717065) -> (7887838,4714657) EMU
$c = $grp.GroupItems.Item("pg57")
$c.Left = 621.0896606445312
$c.Top = 371.23284912109375
# pg58: (7635074,2576424) -> (7674613,2577447) EMU
$c = $grp.GroupItems.Item("pg58")
$c.Left = 604.30029296875
$c.Top = 202.94859313964844
# pg59: (7495520,2237588) -> (7420242,2237553) EMU
$c = $grp.GroupItems.Item("pg59")
$c.Left = 584.2710571289062
$c.Top = 176.18528747558594
# pg60: (8119002,2238840) -> (8199273,2236893) EMU
$c = $grp.GroupItems.Item("pg60")
$c.Left = 645.612060546875
$c.Top = 176.13331604003906
# pg61: (8238462,2238572) -> (7491643,2237824) EMU
$c = $grp.GroupItems.Item("pg61")
$c.Left = 589.8931884765625
$c.Top = 176.2066192626953
# pg62: (8273837,2238763) -> (7945807,2238177) EMU
$c = $grp.GroupItems.Item("pg62")
$c.Left = 625.6541137695312
$c.Top = 176.2344207763672
# pg63: (7864554,2239044) -> (8362116,2238268) EMU
$c = $grp.GroupItems.Item("pg63")
$c.Left = 658.4343872070312
$c.Top = 176.2415771484375
# pg64: (7592628,4838108) -> (7798095,4837745) EMU
$c = $grp.GroupItems.Item("pg64")
$c.Left = 614.0232543945312
$c.Top = 380.9248046875
# pg65: (8060440,4495445) -> (7970369,4496434) EMU
$c = $grp.GroupItems.Item("pg65")
$c.Left = 627.588134765625
$c.Top = 354.0499267578125
# pg66: (7558456,2577780) -> (7606525,2578646) EMU
$c = $grp.GroupItems.Item("pg66")
$c.Left = 598.9390258789062
$c.Top = 203.04299926757812
# pg67: (8036921,3351233) -> (7817374,3352365) EMU
$c = $grp.GroupItems.Item("pg67")
$c.Left = 615.5413208007812
$c.Top = 263.96575927734375
# pg68: (8322289,3204623) -> (8288218,3205487) EMU
$c = $grp.GroupItems.Item("pg68")
$c.Left = 652.6156005859375
$c.Top = 252.4005584716797
# pg69: (7592370,4932255) -> (7698798,4933375) EMU
$c = $grp.GroupItems.Item("pg69")
$c.Left = 606.20458984375
$c.Top = 388.4547424316406
# pg70: (7529331,4554218) -> (8414539,4552730) EMU
$c = $grp.GroupItems.Item("pg70")
$c.Left = 662.5621337890625
$c.Top = 358.4826965332031
# pg71: (8781192,5117024) -> (9343094,5115714) EMU
$c = $grp.GroupItems.Item("pg71")
$c.Left = 735.6766967773438
$c.Top = 402.8121337890625
# pg72: (8353395,2038602) -> (8329471,2038735) EMU
$c = $grp.GroupItems.Item("pg72")
$c.Left = 655.8638916015625
$c.Top = 160.5303192138672
# pg73: (7737763,3127002) -> (7490754,3127725) EMU
$c = $grp.GroupItems.Item("pg73")
$c.Left = 589.8231811523438
$c.Top = 246.27757263183594
# pg74: (8064953,3440302) -> (8122766,3439106) EMU
$c = $grp.GroupItems.Item("pg74")
$c.Left = 639.587890625
$c.Top = 270.7957763671875
# pg75: (8070410,4506508) -> (7708327,4506060) EMU
$c = $grp.GroupItems.Item("pg75")
$c.Left = 606.9548950195312
$c.Top = 354.8078918457031
# pg76: (8126868,2407081) -> (7648895,2409046) EMU
$c = $grp.GroupItems.Item("pg76")
$c.Left = 602.2752075195312
$c.Top = 189.6886749267578
# pg77: (7806830,2276441) -> (7904724,2275635) EMU
$c = $grp.GroupItems.Item("pg77")
$c.Left = 622.4192504882812
$c.Top = 179.18386840820312
# pg78: (8187113,3561658) -> (7841690,3560684) EMU
$c = $grp.GroupItems.Item("pg78")
$c.Left = 617.4559326171875
$c.Top = 280.36883544921875
# pg79: (7643271,4522213) -> (8386479,4522188) EMU
$c = $grp.GroupItems.Item("pg79")
$c.Left = 660.3527221679688
$c.Top = 356.07781982421875
# pg80: (7825012,2160073) -> (8194800,2160372) EMU
$c = $grp.GroupItems.Item("pg80")
$c.Left = 645.2598876953125
$c.Top = 170.1080322265625
# pg81: (7484766,2473511) -> (7763927,2472357) EMU
$c = $grp.GroupItems.Item("pg81")
$c.Left = 611.3328857421875
$c.Top = 194.6737823486328
# pg82: (8335794,2086070) -> (7970422,2085081) EMU
$c = $grp.GroupItems.Item("pg82")
$c.Left = 627.59228515625
$c.Top = 164.1796112060547
# pg83: (7851993,2470163) -> (8240582,2467693) EMU
$c = $grp.GroupItems.Item("pg83")
$c.Left = 648.86474609375
$c.Top = 194.30654907226562
# pg84: (8264776,2034339) -> (8308622,2035629) EMU
$c = $grp.GroupItems.Item("pg84")
$c.Left = 654.2222290039062
$c.Top = 160.28575134277344
# pg85: (8217074,2743155) -> (8306415,2742478) EMU
$c = $grp.GroupItems.Item("pg85")
$c.Left = 654.0484619140625
$c.Top = 215.9431610107422
# pg86: (7817472,1641734) -> (8112328,1641734) EMU
$c = $grp.GroupItems.Item("pg86")
$c.Left = 638.7659912109375
$c.Top = 129.27040100097656
# pg87: (8190751,3233949) -> (8314429,3232743) EMU
$c = $grp.GroupItems.Item("pg87")
$c.Left = 654.6795043945312
$c.Top = 254.5467071533203
# pg88: (7501329,2627866) -> (8343933,2629292) EMU
$c = $grp.GroupItems.Item("pg88")
$c.Left = 657.0026245117188
$c.Top = 207.03086853027344
# pg89: (7817332,2510706) -> (7627328,2510633) EMU
$c = $grp.GroupItems.Item("pg89")
$c.Left = 600.5770263671875
$c.Top = 197.68765258789062
# pg90: (7945936,2726579) -> (7627220,2726045) EMU
$c = $grp.GroupItems.Item("pg90")
$c.Left = 600.5685424804688
$c.Top = 214.6492156982422
# pg91: (7737140,3220454) -> (7441479,3218485) EMU
$c = $grp.GroupItems.Item("pg91")
$c.Left = 585.9432373046875
$c.Top = 253.4240264892578
# pg92: (7896153,2943794) -> (8282460,2944385) EMU
$c = $grp.GroupItems.Item("pg92")
$c.Left = 652.1622314453125
$c.Top = 231.84133911132812
# pg93: (7724708,3094925) -> (7887952,3092118) EMU
$c = $grp.GroupItems.Item("pg93")
$c.Left = 621.0986328125
$c.Top = 243.47386169433594
# pg94: (8094629,3008398) -> (7977073,3008215) EMU
$c = $grp.GroupItems.Item("pg94")
$c.Left = 628.1160278320312
$c.Top = 236.86732482910156
# pg95: (8377501,3293834) -> (8380385,3292379) EMU
$c = $grp.GroupItems.Item("pg95")
$c.Left = 659.8728637695312
$c.Top = 259.2424621582031
# pg96: (8341556,3776349) -> (7543427,3775487) EMU
$c = $grp.GroupItems.Item("pg96")
$c.Left = 593.9706420898438
$c.Top = 297.282470703125
# tx116: (1148183,6810658) -> (1148183,6810354) EMU
$c = $grp.GroupItems.Item("tx116")
$c.Left = 90.40811157226562
$c.Top = 536.2483520507812
# tx117: (1148183,5113465) -> (1148183,5113165) EMU
$c = $grp.GroupItems.Item("tx117")
$c.Left = 90.40811157226562
$c.Top = 402.6114196777344
# tx118: (1148183,3416273) -> (1148183,3415977) EMU
$c = $grp.GroupItems.Item("tx118")
$c.Left = 90.40811157226562
$c.Top = 268.9745788574219
# tx119: (1148183,1719081) -> (1148183,1718789) EMU
$c = $grp.GroupItems.Item("tx119")
$c.Left = 90.40811157226562
$c.Top = 135.3377227783203
# pl120: (1362485,6852349) -> (1362485,6852045) EMU
$c = $grp.GroupItems.Item("pl120")
$c.Left = 107.28228759765625
$c.Top = 539.5311279296875
# pl121: (1362485,5155157) -> (1362485,5154857) EMU
$c = $grp.GroupItems.Item("pl121")
$c.Left = 107.28228759765625
$c.Top = 405.8942565917969
# pl122: (1362485,3457965) -> (1362485,3457669) EMU
$c = $grp.GroupItems.Item("pl122")
$c.Left = 107.28228759765625
$c.Top = 272.2574157714844
# pl123: (1362485,1760772) -> (1362485,1760481) EMU
$c = $grp.GroupItems.Item("pl123")
$c.Left = 107.28228759765625
$c.Top = 138.6205596923828
